$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily record for 2020-04-01 (serial date 43921) as row 44
$ws.Range("A43").Copy()
$ws.Range("A44").PasteSpecial(-4122)
$ws.Range("A44").Value = 43921
$ws.Range("B44").Value = 6777
$ws.Range("C44").Value = 1528
$ws.Range("D44").Value = 81
$ws.Range("E44").Value = 136
$ws.Range("F44").Value = 5249
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0

# Leave the selection where the author ended up after editing
$ws.Range("G43").Select()
